# VerveStacks JPN model update - 2025-08-05 11:56
$wb = $excel.ActiveWorkbook

# --- TS_Defs sheet: Q6 "p,t" -> "t" ---
$wsDefs = $wb.Worksheets.Item("TS_Defs")
$wsDefs.Range("Q6").Value = "t"
$wsDefs.Range("A6").Select()

# --- process map sheet: append two new rows (old_new / */new and old_new / ep*/old) ---
$wsProc = $wb.Worksheets.Item("process map")
$wsProc.Range("A24").Value = "old_new"
$wsProc.Range("B24").Value = "*"
$wsProc.Range("C24").Value = "new"
$wsProc.Range("A25").Value = "old_new"
$wsProc.Range("B25").Value = "ep*"
$wsProc.Range("C25").Value = "old"

# process map becomes the active sheet/selection
$wsProc.Activate()
$wsProc.Range("C26").Select()
